# Finalize BOM: fill in the last (15th) part row -- Red hook-up wire spool --
# and wire up the missing hyperlinks that were left without a clickable link.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New BOM line: row 15 ("Hook-up Wire", Red, 25', Sparkfun PRT-08023) ---
$ws.Range("A15").Value = "Hook-up Wire"
$ws.Range("B15").Value = "Red"
$ws.Range("C15").Value = "25'"
$ws.Range("F15").Value = "Sparkfun"
$ws.Range("G15").Value = "PRT-08023 "
$ws.Range("H15").Value = "https://www.sparkfun.com/products/8023"
$ws.Range("I15").Value = 2.5
$ws.Range("J15").Value = 1
$ws.Range("K15").Formula = "=Table3[[#This Row],[Price]]*Table3[[#This Row],[Quantity]]"

# Remember the existing (already-correct) hyperlink-column styles so we can
# restore them after Hyperlinks.Add(), which otherwise stamps a brand new
# (duplicate) style onto the cell.
$h9Style = $ws.Range("H9").Style
$g15Style = $ws.Range("G15").Style

# --- Missing hyperlinks ---
# H9 already shows the sparkfun URL as text but had no live link.
$ws.Hyperlinks.Add($ws.Range("H9"), "https://www.sparkfun.com/products/8025")
# G15 is the part number of the new row, linking to its product page.
$ws.Hyperlinks.Add($ws.Range("G15"), "https://www.sparkfun.com/products/8023")

$ws.Range("H9").Style = $h9Style
$ws.Range("G15").Style = $g15Style

# Leave the selection where the author last left it.
$ws.Range("K15").Select()
